$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (45171 -> 45172) for every data row (rows 2 through 428).
$ws.Range("C2:C428").Value = 45172
